$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 522, shifting existing rows 522:633 down to 523:634
$ws.Rows.Item(522).Insert()

# Populate the newly inserted row 522 with the new record
$ws.Cells.Item(522, 1).Value = 5
$ws.Cells.Item(522, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(522, 3).Value = "Maule"
$ws.Cells.Item(522, 4).Value = 45244
$ws.Cells.Item(522, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(522, 5).Value = 7
$ws.Cells.Item(522, 6).Value = 100114013
$ws.Cells.Item(522, 7).Value = "Zanahoria"
$ws.Cells.Item(522, 8).Value = "Sin especificar"
$ws.Cells.Item(522, 9).Value = "Primera"
$ws.Cells.Item(522, 10).Value = 600
$ws.Cells.Item(522, 11).Value = 5500
$ws.Cells.Item(522, 12).Value = 5500
$ws.Cells.Item(522, 13).Value = 5500
$ws.Cells.Item(522, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(522, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(522, 16).Value = 275
$ws.Cells.Item(522, 17).Value = 20
$ws.Cells.Item(522, 18).Value = "Hortaliza"
